$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-style existing rows using the Excel built-in cell styles "Good" /
#    "Bad" / "Neutral" (this seeds the new fonts/fills/cellStyleXfs that the
#    target workbook has).
# ---------------------------------------------------------------------------

function Set-RowStyle($row, $styleName, $wrapCol) {
    foreach ($col in @("A","B","C","D")) {
        $cell = $ws.Range($col + $row)
        if ($cell.Value2 -ne $null -or $col -eq $wrapCol -or $true) {
            $cell.Style = $styleName
        }
    }
    if ($wrapCol -ne $null) {
        $ws.Range($wrapCol + $row).WrapText = $true
    }
}

# Row 3: lime -> Bad
Set-RowStyle 3 "Bad" "B"
# Row 4: red -> Good
Set-RowStyle 4 "Good" "B"
# Row 5: stays red -> Bad (no value changes)
Set-RowStyle 5 "Bad" "B"
# Row 6: red -> Neutral (+ value changes, handled below)
Set-RowStyle 6 "Neutral" "B"
# Row 7: lime -> Bad
Set-RowStyle 7 "Bad" "B"
# Row 8: stays red -> Bad
Set-RowStyle 8 "Bad" "B"
# Row 9: yellow -> Bad
Set-RowStyle 9 "Bad" "B"
# Row 10: yellow -> Bad
Set-RowStyle 10 "Bad" "B"
# Row 11: lime -> Bad
Set-RowStyle 11 "Bad" "B"
# Row 12: stays red -> Bad
Set-RowStyle 12 "Bad" "B"
# Row 13: stays red -> Bad
Set-RowStyle 13 "Bad" "B"
# Row 14: stays red -> Bad
Set-RowStyle 14 "Bad" "B"

# D7 has no value but keeps the Bad fill
$ws.Range("D7").Style = "Bad"

# ---------------------------------------------------------------------------
# 2. Row 6 content changes
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "'vM1', 'CA' "
$ws.Range("B6").Value = "vM1', 'CA' "
$ws.Range("D6").Value = "Parece tener periodicidad/solo está bien la corteza"

# ---------------------------------------------------------------------------
# 3. New rows 15-17
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "11361009.abf"
$ws.Range("B15").Value = "Tálamo POM"
$ws.Range("C15").Value = "Si"
$ws.Range("A15:D15").Style = "Good"
$ws.Range("B15").WrapText = $true

$ws.Range("A3").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("A16").Value = "11361091.abf"
$ws.Range("B16").Value = "Tálamo VTM"

$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "16615005.abf"

# ---------------------------------------------------------------------------
# 4. Row heights (Excel auto-increased default row heights across the sheet)
# ---------------------------------------------------------------------------
foreach ($r in 3..15) {
    $ws.Rows($r).RowHeight = 14.5
}

# ---------------------------------------------------------------------------
# 5. Selection
# ---------------------------------------------------------------------------
$ws.Range("C6").Select()
